$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2-12) holds date-like text (e.g. "2025-03-30"). Excel
# would otherwise auto-convert such strings into date serial numbers, so
# force the range to Text format first, assign the values, then restore
# the default "Normal" style so the cells end up as plain shared strings
# using the default style (matching the original workbook layout).
$ws.Range("A2:A12").NumberFormat = "@"

# Header row (row 1)
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

# Data rows 2..12
$ws.Range("A2").Value = "2025-03-30"
$ws.Range("B2").Value = 0.27351930273188851
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("A3").Value = "2025-04-15"
$ws.Range("B3").Value = 0.28444073537437942
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.0047765030929491949
$ws.Range("E3").Value = -0.000022366638060894653
$ws.Range("F3").Value = -0.00086307800864457679
$ws.Range("G3").Value = 0.00050256494321383109
$ws.Range("H3").Value = -0.000086861093565127062
$ws.Range("I3").Value = -0.00037697103016303985
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.00038049635143910798

$ws.Range("A4").Value = "2025-04-30"
$ws.Range("B4").Value = 0.28996684970470421
$ws.Range("C4").Value = -0.0021129358668152995
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.00010172593574537484
$ws.Range("F4").Value = 0.000038058333471570101
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.00014879241856965002
$ws.Range("I4").Value = -0.0020743513521273012
$ws.Range("J4").Value = 0.0003644377745011865
$ws.Range("K4").Value = -0.00071843744835525802

$ws.Range("A5").Value = "2025-05-15"
$ws.Range("B5").Value = 0.299060068643693
$ws.Range("C5").Value = 0.007964037755112818
$ws.Range("D5").Value = -0.0064984566842115265
$ws.Range("E5").Value = -0.00003796853747437684
$ws.Range("F5").Value = 0.00024755226233175604
$ws.Range("G5").Value = -0.0018457040522018122
$ws.Range("H5").Value = -0.000012832643558188193
$ws.Range("I5").Value = -0.00045963658644573245
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.00038762372176409654

$ws.Range("A6").Value = "2025-05-30"
$ws.Range("B6").Value = 0.27022139761345459
$ws.Range("C6").Value = 0.027563494193474362
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -0.000055388575578133195
$ws.Range("F6").Value = -0.000040546798747904394
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.000078462940254440952
$ws.Range("I6").Value = -0.0021853775585362989
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -0.0056317246754117312

$ws.Range("A7").Value = "2025-06-15"
$ws.Range("B7").Value = 0.27268996802050283
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.0019327526525750498
$ws.Range("E7").Value = -0.00049470352567263464
$ws.Range("F7").Value = -0.0043542317416250106
$ws.Range("G7").Value = 0.0011765732653634632
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.00020069058282418058
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.002865914671817571

$ws.Range("A8").Value = "2025-06-30"
$ws.Range("B8").Value = 0.22424356074574131
$ws.Range("C8").Value = -0.055470458870007761
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.00004705853184741944
$ws.Range("F8").Value = -0.000330622086587592
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.000050595149845919418
$ws.Range("I8").Value = 0.0010645514170592124
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.0013660585610660303

$ws.Range("A9").Value = "2025-07-15"
$ws.Range("B9").Value = 0.22316091962422624
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.0034005878414023398
$ws.Range("E9").Value = -0.0034159986345844427
$ws.Range("F9").Value = -0.0046592461356807077
$ws.Range("G9").Value = 0.0013182784748924112
$ws.Range("H9").Value = -0.00015316240076762713
$ws.Range("I9").Value = -0.00033661001488424503
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0.00038561352208421606

$ws.Range("A10").Value = "2025-07-30"
$ws.Range("B10").Value = 0.32640465089051013
$ws.Range("C10").Value = 0.094472038775978887
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -0.00013839053838774564
$ws.Range("F10").Value = -0.00033760059782013413
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = -0.000033708219174563633
$ws.Range("I10").Value = 0.00017031136912818999
$ws.Range("J10").Value = -0.002388063556425647
$ws.Range("K10").Value = 0.0052674881290792896

$ws.Range("A11").Value = "2025-08-15"
$ws.Range("B11").Value = 0.34633560451125067
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = -0.021887310202788182
$ws.Range("E11").Value = 0.0020893345544954326
$ws.Range("F11").Value = 0.0038168772235914984
$ws.Range("G11").Value = 0.0020803197755100704
$ws.Range("H11").Value = 0.00035837715964849218
$ws.Range("I11").Value = 0.00000088567970849626954
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0.0090320524408768477

$ws.Range("A12").Value = "2025-08-30"
$ws.Range("B12").Value = 0.20907350584454676
$ws.Range("C12").Value = -0.061670701625824227
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = -0.00005222843404196606
$ws.Range("F12").Value = 0.000015218656394709551
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0.0000082872376734734251
$ws.Range("I12").Value = -0.0014290518672653819
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = -0.013595073924139611

$ws.Range("A2:A12").Style = "Normal"
